$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the "outbreak" (E) column with "outgroup" for the outgroup rows
# (35, 37, 39) that were previously left blank.
$ws.Cells.Item(35, 5).Value = "outgroup"
$ws.Cells.Item(37, 5).Value = "outgroup"
$ws.Cells.Item(39, 5).Value = "outgroup"

# Fill the "genBankAssembly" (C) and "sha256sumAssembly" (H) columns with "-"
# for every data row (10-39) that doesn't already have a value.
for ($r = 10; $r -le 39; $r++) {
    $ws.Cells.Item($r, 3).Value = "-"
    $ws.Cells.Item($r, 8).Value = "-"
}

# Row 14's genBankAssembly cell had stray bold formatting left over from the
# header style; normalize it to match the rest of the column.
$ws.Cells.Item(14, 3).Font.Bold = $false

# Update the saved view/selection state to match the authored edit.
$ws.Range("H16").Select()
